$wb = $excel.ActiveWorkbook

# --- Work on the "2023 League Hitting" sheet ---
$ws = $wb.Worksheets.Item("2023 League Hitting")

# Rename the "BA" header (column R) to "AVG" - same underlying data, just a
# relabeled header.
$ws.Range("R1").Value = "AVG"

# Add the new "wRC+" column in AE, right after the existing WAR column (AD).
# Give the header the same look as the other header cells (copy format from
# the WAR header) and the data cells the same look as the WAR data column.
$ws.Range("AD1").Copy()
$ws.Range("AE1").PasteSpecial(-4122) | Out-Null
$ws.Range("AE1").Value = "wRC+"

$ws.Range("AD2:AD31").Copy()
$ws.Range("AE2:AE31").PasteSpecial(-4122) | Out-Null

$wrcValues = @(96,125,105,99,106,83,97,91,80,90,112,86,101,116,92,92,108,100,92,88,106,90,106,108,91,103,117,116,107,93)
for ($i = 0; $i -lt $wrcValues.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 31)
    $cell.Value = $wrcValues[$i]
}

# League-average row: AVERAGE formula, right-aligned, smaller non-bold font,
# one-decimal-place number format (matches the WAR column's averaging style).
$avgCell = $ws.Range("AE32")
$avgCell.Formula = "=AVERAGE(AE2:AE31)"
$avgCell.NumberFormat = "0.0"
$avgCell.Font.Size = 10
$avgCell.Font.Bold = $false
$avgCell.HorizontalAlignment = -4152

# Make the Hitting sheet the active tab and move the selection to AF31, as in
# the target workbook.
$ws.Activate()
$ws.Range("AF31").Select() | Out-Null
